$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "ID" header in A1 (creates a new shared string entry)
$ws.Range("A1").Value = "ID"

# Remove the extra styled-but-empty columns H:P (rows 1-13 had s="4" placeholders)
$ws.Range("H1:P13").Delete()

# Remove the now-unused trailing rows 14-18 (they only held s="4" placeholders)
$ws.Rows("14:18").Delete()

# Move/restore the active selection to A2, matching the saved view state
$ws.Range("A2").Select()
